$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells P1 and Q1, matching the style/format of O1 (bold/border/centered header style)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15

# For each data row 2-25: swap I<->K and M<->O, then append P=2, Q=2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # I
    $kVal = $ws.Cells.Item($r, 11).Value2  # K
    $ws.Cells.Item($r, 9).Value2 = $kVal
    $ws.Cells.Item($r, 11).Value2 = $iVal

    $mVal = $ws.Cells.Item($r, 13).Value2  # M
    $oVal = $ws.Cells.Item($r, 15).Value2  # O
    $ws.Cells.Item($r, 13).Value2 = $oVal
    $ws.Cells.Item($r, 15).Value2 = $mVal

    $ws.Cells.Item($r, 16).Value2 = 2  # P
    $ws.Cells.Item($r, 17).Value2 = 2  # Q
}
